# Updates currency/price snapshot values (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the 8 crafting-job leve-profit sheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR) to reflect freshly pulled market-board data.
#
# LeveProfitNQ (col M) is only present when LevePriceNQ (col K) > 0, and
# LeveProfitHQ (col N) is only present when LevePriceHQ (col L) > 0 -- a
# couple of rows below gain/lose one of those two cells entirely as their
# K/L values cross zero.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1083.25
$ws.Range("J17").Value = 1166.5714
$ws.Range("L17").Value = 3499.7142
$ws.Range("N17").Value = -3835.7142

$ws.Range("H33").Value = 234.75
$ws.Range("I33").Value = 234.75
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 234.75
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -5.75
$ws.Range("N33").ClearContents()

$ws.Range("H80").Value = 674.8461
$ws.Range("I80").Value = 441
$ws.Range("J80").Value = 1049
$ws.Range("K80").Value = 1323
$ws.Range("L80").Value = 3147
$ws.Range("M80").Value = -325
$ws.Range("N80").Value = -5143

$ws.Range("H83").Value = 674.8461
$ws.Range("I83").Value = 441
$ws.Range("J83").Value = 1049
$ws.Range("K83").Value = 3969
$ws.Range("L83").Value = 9441
$ws.Range("M83").Value = 1023
$ws.Range("N83").Value = -19425

$ws.Range("H103").Value = 3851.52
$ws.Range("I103").Value = 2160.2222
$ws.Range("J103").Value = 4802.875
$ws.Range("K103").Value = 6480.6666
$ws.Range("L103").Value = 14408.625
$ws.Range("M103").Value = -5894.6666
$ws.Range("N103").Value = -15580.625

$ws.Range("H137").Value = 2734.1052
$ws.Range("I137").Value = 1954.4166
$ws.Range("J137").Value = 4070.7144
$ws.Range("K137").Value = 5863.2498
$ws.Range("L137").Value = 12212.1432
$ws.Range("M137").Value = -3313.2498
$ws.Range("N137").Value = -17312.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 846.3182
$ws.Range("I2").Value = 730.7857
$ws.Range("J2").Value = 1048.5
$ws.Range("K2").Value = 730.7857
$ws.Range("L2").Value = 1048.5
$ws.Range("M2").Value = -617.7857
$ws.Range("N2").Value = -1274.5

$ws.Range("H38").Value = 3649.6667
$ws.Range("I38").Value = 3649.6667
$ws.Range("K38").Value = 3649.6667
$ws.Range("M38").Value = -3182.6667

$ws.Range("H63").Value = 7099.222
$ws.Range("I63").Value = 1224.5
$ws.Range("J63").Value = 11799
$ws.Range("K63").Value = 1224.5
$ws.Range("L63").Value = 11799
$ws.Range("M63").Value = -538.5
$ws.Range("N63").Value = -13171

$ws.Range("H66").Value = 7099.222
$ws.Range("I66").Value = 1224.5
$ws.Range("J66").Value = 11799
$ws.Range("K66").Value = 6122.5
$ws.Range("L66").Value = 58995
$ws.Range("M66").Value = -2690.5
$ws.Range("N66").Value = -65859

$ws.Range("H116").Value = 846.3182
$ws.Range("I116").Value = 730.7857
$ws.Range("J116").Value = 1048.5
$ws.Range("K116").Value = 730.7857
$ws.Range("L116").Value = 1048.5
$ws.Range("M116").Value = 1563.2143
$ws.Range("N116").Value = -5636.5

$ws.Range("H132").Value = 1767.9375
$ws.Range("I132").Value = 1619.1333
$ws.Range("K132").Value = 4857.3999
$ws.Range("M132").Value = -2327.3999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 846.3182
$ws.Range("I3").Value = 730.7857
$ws.Range("J3").Value = 1048.5
$ws.Range("K3").Value = 730.7857
$ws.Range("L3").Value = 1048.5
$ws.Range("M3").Value = -616.7857
$ws.Range("N3").Value = -1276.5

$ws.Range("H134").Value = 3420.375
$ws.Range("I134").Value = 3593.8333
$ws.Range("J134").Value = 2900
$ws.Range("K134").Value = 10781.4999
$ws.Range("L134").Value = 8700
$ws.Range("M134").Value = -8246.499899999999
$ws.Range("N134").Value = -13770

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H38").Value = 4686.727
$ws.Range("I38").Value = 6380.857
$ws.Range("J38").Value = 1722
$ws.Range("K38").Value = 6380.857
$ws.Range("L38").Value = 1722
$ws.Range("M38").Value = -6003.857
$ws.Range("N38").Value = -2476

$ws.Range("H46").Value = 4686.727
$ws.Range("I46").Value = 6380.857
$ws.Range("J46").Value = 1722
$ws.Range("K46").Value = 6380.857
$ws.Range("L46").Value = 1722
$ws.Range("M46").Value = -6169.857
$ws.Range("N46").Value = -2144

$ws.Range("H58").Value = 3298.8
$ws.Range("I58").Value = 1850
$ws.Range("K58").Value = 1850
$ws.Range("M58").Value = -1647

$ws.Range("H136").Value = 3298.8
$ws.Range("I136").Value = 1850
$ws.Range("K136").Value = 5550
$ws.Range("M136").Value = -3000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 126.833336
$ws.Range("I6").Value = 47.88889
$ws.Range("K6").Value = 143.66667
$ws.Range("M6").Value = -30.66667000000001

$ws.Range("H40").Value = 282.44446
$ws.Range("J40").Value = 361.57144
$ws.Range("L40").Value = 1446.28576
$ws.Range("N40").Value = -1584.28576

$ws.Range("H113").Value = 642.35
$ws.Range("I113").Value = 542.1667
$ws.Range("J113").Value = 792.625
$ws.Range("K113").Value = 1626.5001
$ws.Range("L113").Value = 2377.875
$ws.Range("M113").Value = 543.4999
$ws.Range("N113").Value = -6717.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2670.4
$ws.Range("I132").Value = 2418.2307
$ws.Range("K132").Value = 7254.6921
$ws.Range("M132").Value = -4724.6921

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 7033.3335
$ws.Range("J46").Value = 7033.3335
$ws.Range("L46").Value = 7033.3335
$ws.Range("N46").Value = -7409.3335

$ws.Range("H82").Value = 7249
$ws.Range("I82").Value = 6247.5
$ws.Range("J82").Value = 7499.375
$ws.Range("K82").Value = 6247.5
$ws.Range("L82").Value = 7499.375
$ws.Range("M82").Value = -5886.5
$ws.Range("N82").Value = -8221.375

$ws.Range("H85").Value = 7249
$ws.Range("I85").Value = 6247.5
$ws.Range("J85").Value = 7499.375
$ws.Range("K85").Value = 6247.5
$ws.Range("L85").Value = 7499.375
$ws.Range("M85").Value = -4999.5
$ws.Range("N85").Value = -9995.375

$ws.Range("H96").Value = 55197
$ws.Range("J96").Value = 55197
$ws.Range("L96").Value = 55197
$ws.Range("N96").Value = -60689

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8544.333000000001
$ws.Range("I81").Value = 6000
$ws.Range("J81").Value = 11724.75
$ws.Range("K81").Value = 12000
$ws.Range("L81").Value = 23449.5
$ws.Range("M81").Value = -10939
$ws.Range("N81").Value = -25571.5

$ws.Range("H84").Value = 8544.333000000001
$ws.Range("I84").Value = 6000
$ws.Range("J84").Value = 11724.75
$ws.Range("K84").Value = 60000
$ws.Range("L84").Value = 117247.5
$ws.Range("M84").Value = -54696
$ws.Range("N84").Value = -127855.5

$ws.Range("H126").Value = 6000.1113
$ws.Range("I126").Value = 3371
$ws.Range("K126").Value = 10113
$ws.Range("M126").Value = -7643

$ws.Range("H136").Value = 2983.2856
$ws.Range("I136").Value = 2261.0908
$ws.Range("K136").Value = 6783.2724
$ws.Range("M136").Value = -4233.2724
